# Update the cryptos worksheet with refreshed price/volume data.
# A new coin ("OKB") entered the top-50 ranking, which shifts several
# existing rows down by one position and drops "Cronos" off the list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    # Force a text number format so numeric-looking strings (e.g. "1.001")
    # are not silently coerced into floating point numbers by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    # Restore the default/normal style so we do not leave a stray text
    # number-format style attached to the cell.
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 "29.976.85"
Set-TextCell 2 5 "  -0.82%  "
# Row 3
Set-TextCell 3 4 "1.919.03"
Set-TextCell 3 5 "  +1.21%  "
# Row 4
Set-TextCell 4 4 "1.001"
Set-TextCell 4 5 "  -0.01%  "
# Row 5
Set-TextCell 5 4 "320.11"
Set-TextCell 5 5 "  -1.53%  "
# Row 6
Set-TextCell 6 4 "1.000"
Set-TextCell 6 5 "  +0.00%  "
# Row 7
Set-TextCell 7 4 "0.5060"
Set-TextCell 7 5 "  -2.36%  "
# Row 8
Set-TextCell 8 4 "0.4045"
Set-TextCell 8 5 "  +0.97%  "
# Row 9
Set-TextCell 9 4 "0.08354"
Set-TextCell 9 5 "  -0.44%  "
# Row 10
Set-TextCell 10 4 "1.105"
Set-TextCell 10 5 "  -0.85%  "
# Row 11
Set-TextCell 11 2 "OKB"
Set-TextCell 11 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell 11 4 "42.03"
Set-TextCell 11 5 "  -1.59%  "
# Row 12
Set-TextCell 12 2 "Solana"
Set-TextCell 12 3 "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell 12 4 "23.91"
Set-TextCell 12 5 "  +2.44%  "
# Row 13
Set-TextCell 13 2 "WrappedEther"
Set-TextCell 13 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell 13 4 "1.919.03"
Set-TextCell 13 5 "  +1.37%  "
# Row 14
Set-TextCell 14 2 "Polkadot"
Set-TextCell 14 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell 14 4 "6.406"
Set-TextCell 14 5 "  -0.46%  "
# Row 15
Set-TextCell 15 2 "Chainlink"
Set-TextCell 15 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell 15 4 "7.239"
Set-TextCell 15 5 "  -1.23%  "
# Row 16
Set-TextCell 16 2 "BinanceUSD"
Set-TextCell 16 3 "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell 16 4 "1.000"
Set-TextCell 16 5 "  -0.08%  "
# Row 17
Set-TextCell 17 2 "Litecoin"
Set-TextCell 17 3 "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell 17 4 "92.31"
Set-TextCell 17 5 "  -2.11%  "
# Row 18
Set-TextCell 18 2 "ShibaInu"
Set-TextCell 18 3 "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell 18 4 "0.00001099"
Set-TextCell 18 5 "  -0.99%  "
# Row 19
Set-TextCell 19 2 "TRON"
Set-TextCell 19 3 "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell 19 4 "0.06511"
Set-TextCell 19 5 "  -2.05%  "
# Row 20
Set-TextCell 20 2 "Avalanche"
Set-TextCell 20 3 "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell 20 4 "18.24"
Set-TextCell 20 5 "  +0.14%  "
# Row 21
Set-TextCell 21 2 "Dai"
Set-TextCell 21 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell 21 4 "1.001"
Set-TextCell 21 5 "  +0.05%  "
# Row 22
Set-TextCell 22 2 "Uniswap"
Set-TextCell 22 3 "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell 22 4 "5.950"
Set-TextCell 22 5 "  +0.00%  "
# Row 23
Set-TextCell 23 2 "WrappedBTC"
Set-TextCell 23 3 "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell 23 4 "30.003.61"
Set-TextCell 23 5 "  -0.70%  "
# Row 24
Set-TextCell 24 2 "Cosmos"
Set-TextCell 24 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell 24 4 "11.34"
Set-TextCell 24 5 "  +0.42%  "
# Row 25
Set-TextCell 25 2 "Toncoin"
Set-TextCell 25 3 "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell 25 4 "2.196"
Set-TextCell 25 5 "  -1.75%  "
# Row 26
Set-TextCell 26 4 "22.16"
Set-TextCell 26 5 "  +2.56%  "
# Row 27
Set-TextCell 27 2 "WrappedliquidstakedEther2.0"
Set-TextCell 27 3 "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell 27 4 "2.137.95"
Set-TextCell 27 5 "  +1.32%  "
# Row 28
Set-TextCell 28 2 "Monero"
Set-TextCell 28 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell 28 4 "162.46"
Set-TextCell 28 5 "  +0.42%  "
# Row 29
Set-TextCell 29 2 "LidoDAOToken"
Set-TextCell 29 3 "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell 29 4 "2.330"
Set-TextCell 29 5 "  -0.89%  "
# Row 30
Set-TextCell 30 2 "BitcoinCash"
Set-TextCell 30 3 "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell 30 4 "129.20"
Set-TextCell 30 5 "  -0.42%  "
# Row 31
Set-TextCell 31 2 "ImmutableX"
Set-TextCell 31 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell 31 4 "1.130"
Set-TextCell 31 5 "  +3.40%  "
# Row 32
Set-TextCell 32 2 "Stellar"
Set-TextCell 32 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 32 4 "0.1039"
Set-TextCell 32 5 "  -1.40%  "
# Row 33
Set-TextCell 33 2 "Filecoin"
Set-TextCell 33 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 33 4 "5.969"
Set-TextCell 33 5 "  -1.90%  "
# Row 34
Set-TextCell 34 2 "HuobiToken"
Set-TextCell 34 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell 34 4 "3.818"
Set-TextCell 34 5 "  +2.04%  "
# Row 35
Set-TextCell 35 2 "VeChain"
Set-TextCell 35 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell 35 4 "0.02453"
Set-TextCell 35 5 "  -1.69%  "
# Row 36
Set-TextCell 36 2 "InternetComputer(DFINITY)"
Set-TextCell 36 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell 36 4 "5.399"
Set-TextCell 36 5 "  +2.68%  "
# Row 37
Set-TextCell 37 2 "Hedera"
Set-TextCell 37 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 37 4 "0.06421"
Set-TextCell 37 5 "  -1.95%  "
# Row 38
Set-TextCell 38 2 "Algorand"
Set-TextCell 38 3 "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell 38 4 "0.2157"
Set-TextCell 38 5 "  -1.75%  "
# Row 39
Set-TextCell 39 4 "8.754"
Set-TextCell 39 5 "  +0.19%  "
# Row 40
Set-TextCell 40 4 "0.6493"
Set-TextCell 40 5 "  -0.11%  "
# Row 41
Set-TextCell 41 2 "ARBITRUM"
Set-TextCell 41 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 41 4 "1.194"
Set-TextCell 41 5 "  -2.19%  "
# Row 42
Set-TextCell 42 2 "Aptos"
Set-TextCell 42 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell 42 4 "11.38"
Set-TextCell 42 5 "  -3.44%  "
# Row 43
Set-TextCell 43 2 "TrustWalletToken"
Set-TextCell 43 3 "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell 43 4 "1.215"
Set-TextCell 43 5 "  -0.87%  "
# Row 44
Set-TextCell 44 2 "NEARProtocol"
Set-TextCell 44 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell 44 4 "2.221"
Set-TextCell 44 5 "  +8.21%  "
# Row 45
Set-TextCell 45 2 "EnergySwap"
Set-TextCell 45 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell 45 4 "13.44"
Set-TextCell 45 5 "  +1.22%  "
# Row 46
Set-TextCell 46 2 "Decentraland"
Set-TextCell 46 3 "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell 46 4 "0.6053"
Set-TextCell 46 5 "  -0.65%  "
# Row 47
Set-TextCell 47 2 "PancakeSwap"
Set-TextCell 47 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell 47 4 "3.640"
Set-TextCell 47 5 "  -1.60%  "
# Row 48
Set-TextCell 48 5 "  -1.84%  "
# Row 49
Set-TextCell 49 2 "EOS"
Set-TextCell 49 3 "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextCell 49 4 "1.210"
Set-TextCell 49 5 "  -2.20%  "
# Row 50
Set-TextCell 50 2 "Aave"
Set-TextCell 50 3 "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell 50 4 "79.03"
Set-TextCell 50 5 "  +0.11%  "
# Row 51
Set-TextCell 51 2 "WEMIXTOKEN"
Set-TextCell 51 3 "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell 51 4 "1.130"
Set-TextCell 51 5 "  -2.88%  "
